$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5848
$ws.Range("K3").Value = 6024
$ws.Range("C4").Value = 1539
$ws.Range("F4").Value = 1578
$ws.Range("K4").Value = 1252
$ws.Range("K5").Value = 426
$ws.Range("K6").Value = 6630
$ws.Range("C7").Value = 22626
$ws.Range("F7").Value = 19400
$ws.Range("K7").Value = 20180

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 584
$ws.Range("K8").Value = 1336
$ws.Range("K9").Value = 88
$ws.Range("K11").Value = 382
$ws.Range("K15").Value = 208
$ws.Range("K17").Value = 39
$ws.Range("K18").Value = 133
$ws.Range("K19").Value = 583
$ws.Range("K20").Value = 476
$ws.Range("K23").Value = 208
$ws.Range("K25").Value = 96
$ws.Range("K26").Value = 25
$ws.Range("K27").Value = 188
$ws.Range("K29").Value = 1103
$ws.Range("K31").Value = 225
$ws.Range("K33").Value = 873
$ws.Range("K34").Value = 113
$ws.Range("K35").Value = 32
$ws.Range("K37").Value = 681
$ws.Range("K41").Value = 139
$ws.Range("K42").Value = 748
$ws.Range("K43").Value = 173
$ws.Range("K47").Value = 140
$ws.Range("K52").Value = 531
$ws.Range("K53").Value = 256
$ws.Range("K54").Value = 392
$ws.Range("K55").Value = 224
$ws.Range("K57").Value = 77
$ws.Range("K58").Value = 15
$ws.Range("C63").Value = 229
$ws.Range("F63").Value = 157
$ws.Range("K63").Value = 57
$ws.Range("K65").Value = 467
$ws.Range("K67").Value = 787
$ws.Range("K72").Value = 95
$ws.Range("K73").Value = 178
$ws.Range("K76").Value = 273
$ws.Range("K77").Value = 141
$ws.Range("K78").Value = 227
$ws.Range("K83").Value = 450
$ws.Range("K88").Value = 217
$ws.Range("K89").Value = 296
$ws.Range("K91").Value = 229
$ws.Range("K94").Value = 270
$ws.Range("K95").Value = 340
$ws.Range("K96").Value = 213
$ws.Range("K98").Value = 93
$ws.Range("K99").Value = 333
$ws.Range("C101").Value = 22626
$ws.Range("F101").Value = 19400
$ws.Range("K101").Value = 20180

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 189
$ws.Range("K4").Value = 21
$ws.Range("K6").Value = 156
$ws.Range("K7").Value = 584

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 135
$ws.Range("K7").Value = 382

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 35
$ws.Range("K7").Value = 296

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 188
$ws.Range("K7").Value = 531

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 370
$ws.Range("K3").Value = 407
$ws.Range("K4").Value = 74
$ws.Range("K6").Value = 447
$ws.Range("K7").Value = 1336

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 156
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 450

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 322
$ws.Range("K7").Value = 873

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 115
$ws.Range("K3").Value = 120
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 199
$ws.Range("K3").Value = 225
$ws.Range("K7").Value = 681

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 114
$ws.Range("K7").Value = 467

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 137
$ws.Range("K7").Value = 333

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K3").Value = 60
$ws.Range("K5").Value = 7
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 281
$ws.Range("K6").Value = 223
$ws.Range("K7").Value = 787

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 96
$ws.Range("K6").Value = 210
$ws.Range("K7").Value = 392

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 314
$ws.Range("K3").Value = 396
$ws.Range("K4").Value = 53
$ws.Range("K6").Value = 312
$ws.Range("K7").Value = 1103

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 174
$ws.Range("K3").Value = 178
$ws.Range("K4").Value = 26
$ws.Range("K6").Value = 187
$ws.Range("K7").Value = 583

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 273

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 202
$ws.Range("K3").Value = 230
$ws.Range("K6").Value = 279
$ws.Range("K7").Value = 748

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 69
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 158
$ws.Range("K3").Value = 155
$ws.Range("K7").Value = 476

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 270

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 40
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 74
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K3").Value = 19
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 60
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 15
